# Append a new day (2025-12-30) to the GSC export "Chart" sheet.
# The new row mirrors the previous day's row: Invalid = 0, Valid = 28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$newRow = 87

# Force column A to be stored as plain text (matching the existing date
# strings in the sheet) instead of being auto-converted to a date serial.
# A leading apostrophe makes Excel treat the value as a text literal; the
# apostrophe itself is a quote-prefix qualifier and is not stored in the
# cell's value. ClearFormats() then drops the quote-prefix formatting flag
# that the apostrophe entry leaves behind, returning the cell to the same
# (default) style as the rest of the column.
$ws.Cells.Item($newRow, 1).Value = "'2025-12-30"
$ws.Cells.Item($newRow, 1).ClearFormats()
$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 28
